{"js": "// Append \"  \" (two trailing spaces) to the first paragraph's existing\n// text, then append the red \"(This is a change \u2013 Version for main\n// branch)\" annotation as three separate colored runs, matching the\n// author's incremental edit captured in the diff.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// 1) Pad the original sentence with two trailing spaces (stays black).\nfirstParagraph.insertText(\"  \", Word.InsertLocation.end);\nawait context.sync();\n\n// 2) Append the red annotation in three runs (as in the source diff).\nconst run1 = firstParagraph.insertText(\"(This is a change \\u2013 Ve\", Word.InsertLocation.end);\nrun1.font.color = \"#FF0000\";\nawait context.sync();\n\nconst run2 = firstParagraph.insertText(\"rsion for main branch\", Word.InsertLocation.end);\nrun2.font.color = \"#FF0000\";\nawait context.sync();\n\nconst run3 = firstParagraph.insertText(\")\", Word.InsertLocation.end);\nrun3.font.color = \"#FF0000\";\nawait context.sync();\n", "ps1": "# Append \"  \" (two trailing spaces) to the first paragraph's existing\n# text, then append the red \"(This is a change \u2013 Version for main\n# branch)\" annotation as three separate colored runs, matching the\n# author's incremental edit captured in the diff.\n\n$d = $word.ActiveDocument\n$para = $d.Paragraphs(1).Range\n\n# 1) Pad the original sentence with two trailing spaces (stays black).\n#    InsertAfter on the paragraph's own Range is paragraph-mark aware,\n#    so the new text lands before the pilcrow and $para grows to match.\n$para.InsertAfter(\"  \")\n\n# 2) Append the red annotation as three runs, coloring each newly\n#    inserted span right after creating it. $para.End always points\n#    one past the paragraph mark, so the actual insertion point is\n#    $para.End - 1.\n$start = $para.End - 1\n$para.InsertAfter(\"(This is a change \" + [char]0x2013 + \" Ve\")\n$d.Range($start, $para.End - 1).Font.Color = 255\n\n$start = $para.End - 1\n$para.InsertAfter(\"rsion for main branch\")\n$d.Range($start, $para.End - 1).Font.Color = 255\n\n$start = $para.End - 1\n$para.InsertAfter(\")\")\n$d.Range($start, $para.End - 1).Font.Color = 255\n"}
